$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 173916
$ws.Range("C4").Value = 163916
$ws.Range("C5").Value = 10000
$ws.Range("C6").Value = 900
$ws.Range("C7").Value = 5.75
$ws.Range("C8").Value = 64.36
